$d = $word.ActiveDocument

$found = $d.Content.Find.Execute('Geological - geotechnical constraints in environmental studies', $true, $false, $false, $false, $false, $true, 1, $false, 'Environmental Geotechnics', 2)
Write-Output "Replace 0: $found"

$found = $d.Content.Find.Execute('Ativação: 01/01/2020', $true, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2022', 2)
Write-Output "Replace 1: $found"

$found = $d.Content.Find.Execute('Compreensão dos principais condicionantes atuantes nos processos de movimentos de massa, erosão, contaminação e em áreas que serão utilizadas como aterros.', $true, $false, $false, $false, $false, $true, 1, $false, 'Propiciar o entendimento sobre os conceitos básicos, teoria, metodologias de análise e práticas de engenharia em geotecnia ambiental.', 2)
Write-Output "Replace 2: $found"

$found = $d.Content.Find.Execute('Understanding of the main conditioning factors in the processes of mass movements, erosion, contamination and in areas that will be used as landfills.', $true, $false, $false, $false, $false, $true, 1, $false, 'Provide knowledge of the basic concepts, theory,  analysis methodologies and engineering practices in environmental geotechnics..', 2)
Write-Output "Replace 3: $found"

$found = $d.Content.Find.Execute('1Introdução – conceitos (condicionantes, agentes deflagradores, atributos, etc);2Materiais terrestres; Abordagem de algumas propriedades e parâmetros dos materiais terrestres que são importantes em estudos de condicionantes geológicos-geotécnicos;3Movimento de massa – análise dos condicionantes e dos processos deflagradores. Estudo de caso;4Erosão – análise dos condicionantes e dos processos deflagradores. Estudo de caso;5Contaminação – análise dos condicionantes. Estudo de caso;5.1Contaminantes orgânicos5.2Contaminantes inorgânicos6Aterros – análise dos condicionantes. Estudo de caso;7Atividades humanas, condicionantes geológico-geotécnicos e estudos ambientais.', $true, $false, $false, $false, $false, $true, 1, $false, 'Principais fatores condicionantes atuantes em processos como erosão, movimento gravitacional de massa e contaminação de materiais geológicos e água subterrânea. Técnicas de investigação, monitoramento, contenção e recuperação em geotecnia ambiental.', 2)
Write-Output "Replace 4: $found"

$found = $d.Content.Find.Execute('1Introduction - concepts (aspects, triggers, attributes, etc);2Terrestrial materials; Approach of some properties and parameters of terrestrial materials that are important in studies of geological-geotechnical aspects;3Mass movement - analysis of the geological geotechnical aspects and the triggering processes. Case study;4Erosion - analysis of the geological geotechnical aspects and the triggering processes. Case study;5Contamination - analysis of the geological geotechnical aspects. Case study;5.1Organic contaminants5.2Inorganic contaminants6Landfills - analysis of the geological geotechnical aspects. Case study;7Human activities, geological-geotechnical conditions and environmental studies.', $true, $false, $false, $false, $false, $true, 1, $false, 'Main aspects acting in processes such as erosion, gravitational mass movements and contamination of geological materials and groundwater. Investigation, monitoring, containment and recuperation techniques in environmental geotechnics.', 2)
Write-Output "Replace 5: $found"

$found = $d.Content.Find.Execute('Introdução à disciplina: abordagem dos principais conceitos. Definição: condicionantes e agentes deflagradores. Materiais terrestres. Granulometria do solo, porosidade, permeabilidade, teor de matéria orgânica, mineralogia, materiais carbonáticos, declividade, zona de recarga, nível d’água. Ação destes atributos como fatores condicionantes em áreas degradadas. Movimento de massa: principais tipos de movimentos de massa. Condicionantes geológico-geotécnicos relacionados com os movimentos de massa. Erosão: erosão hídrica – laminar e linear. Implicações com uso e ocupação humano. Condicionantes geológico-geotécnicos relacionados com a erosão hídrica. Contaminação: Conceitos de contaminação x poluição, fontes de contaminação, tipos de contaminantes, valores orientadores (background, CONAMA e Cetesb). Condicionantes geológico-geotécnicos relacionados com a contaminação de contaminantes orgânicos. Condicionantes geológico-geotécnicos relacionados com a contaminação de contaminantes inorgânicos. Aterros: escolha da área para instalação de um aterro de resíduos sólidos. Condicionantes geológico-geotécnicos relacionados com a escolha de área para aterro. Aplicação do conhecimento dos condicionantes geológico-geotécnicos em estudos ambientais.', $true, $false, $false, $false, $false, $true, 1, $false, 'Aspectos conceituais e abrangência do tema. Investigação e monitoramento geoambiental. Mapas geoambientais e geotécnicos. Movimentos gravitacionais de massa: conceitos, tipos e mecanismos. Análise de estabilidade de taludes e técnicas de estabilização de encostas naturais. Processos erosivos. Estimativas de perda de solo e técnicas de mitigação de processos erosivos antrópicos. Contaminação de materiais geológicos e água subterrânea. Transporte de contaminantes na zona saturada e não saturada do solo. Barragens de rejeitos e aterros de resíduos. Geossintéticos', 2)
Write-Output "Replace 6: $found"

$found = $d.Content.Find.Execute('Introduction to the discipline: approach of the main concepts. Definition: aspects and trigger agents. Terrestrial materials. Soil granulometry, porosity, permeability, organic matter content, mineralogy, carbonate materials, slope, recharge zone, water level. Action of these attributes as conditioning factors in degraded areas. Mass movement: main types of mass movements. Geological-geotechnical factors related to mass movements. Erosion: water erosion - laminar and linear. Implications with human use and occupation. Geological-geotechnical factors related to water erosion. Contamination: Concepts of contamination x pollution, sources of contamination, types of contaminants, guiding values (background, CONAMA and Cetesb). Geological-geotechnical factors related to the contamination of organic contaminants. Geological-geotechnical factors related to the contamination of inorganic contaminants. Landfills: choice of area for installation of a solid waste landfill. Geological-geotechnical factors related to the choice of landfill area. Application of geological-geotechnical aspects in environmental studies.', $true, $false, $false, $false, $false, $true, 1, $false, 'Conceptual aspects and scope of the theme. Geoenvironmental investigation and monitoring. Geoenvironmental and geotechnical maps. Gravitational mass movements: concepts, types and mechanics. Slope stability analysis and natural slope stabilization techniques. Erosive processes. Estimates of soil loss and mitigation techniques for anthropic erosive processes. Contamination of geological materials and groundwater. Transport of contaminants in saturated and unsaturated zone of the soil. Tailings dams and waste landfills. Geosynthetics.', 2)
Write-Output "Replace 7: $found"

$found = $d.Content.Find.Execute('Aulas expositivas teóricas e práticas; Atividades discentes: pesquisas bibliográficas e seminários', $true, $false, $false, $false, $false, $true, 1, $false, 'Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas.', 2)
Write-Output "Replace 8: $found"

$found = $d.Content.Find.Execute('Média ponderada de exercícios e provas.', $true, $false, $false, $false, $false, $true, 1, $false, 'Média ponderada de provas  e atividades.', 2)
Write-Output "Replace 9: $found"

$found = $d.Content.Find.Execute('A nota final (MF) do aluno que realizou prova de recuperação dependerá da média do semestre (MS) e da média da prova de recuperação (MR), como segue:MF=5 se 5 ≤MR ≤ (10 - MS); MF = (MS + MR) / 2 se MR > (10 – MS)MF = MS se MR < 5.', $true, $false, $false, $false, $false, $true, 1, $false, '1 (uma) prova escrita', 2)
Write-Output "Replace 10: $found"

$found = $d.Content.Find.Execute('•BARROW, C. J. Land degradation, 1991. •BITAR, O. Curso de geologia aplicada ao meio ambiente. ABGE, 1995. •BOSCOV, M. E. Geotecnia ambiental. Oficina de Textos, 2008. •CALIJURI, M. C.; CUNHA, D. G. F. Engenharia ambiental: Conceitos, tecnologia e gestão. Elsevier, 2012.•GUERRA, A. J. T.; CUNHA, S. B. Geomorfologia e meio ambiente, 1996.•OLIVEIRA, A. M. S.; BRITO, S. N. A. Geologia de engenharia. ABGE, 1998. •PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra, 2008. •SHARMA, H. D.; REDDY, K. R. Geoenvironmental engineering, 2004. •ZUQUETTE, L. V. Geotecnia ambiental. Elsevier, 2015.•YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. 2001.', $true, $false, $false, $false, $false, $true, 1, $false, 'Bibliografia básica:CHIOSSI, n. Geologia de Engenharia. Oficina de Textos: São Paulo, 1979. 427p.Bibliografia complementarBOSCOV, M. E. Geotecnia ambiental. Oficina de Textos: São Paulo, 2008. 248p.ZUQUETTE, L. V. Geotecnia ambiental. Elsevier: Rio de Janeiro, 2015. 432p.OLIVEIRA, A. M. S.; JERÔNIMO, J. Geologia de Engenharia e Ambiental, ABGE, 2018. 912 p.BARROW, C. J. Land degradation. Cambridge University Press: New York, 1991. 316p.KHALED, S. B. D. Fundamentos de Engenharia Geotécnica. Cengage Learning. 2020, 630p.SHARMA, H. D.; REDDY, K. R. Geoenvironmental engineering, Wiley, 2004. 992p. YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p.DANIEL, D.E.  Geotechnical practice for waste disposal. Chapman & Hall, 1993. 693p.ROWE, R. K.; QUIGLEY, R.M.; BOOKER, J. R. Clayey barrier systems for waste disposal facilities. E & FN SPON, 2019. 600p.', 2)
Write-Output "Replace 11: $found"
